# Deleted post rule check from Excel, Template and variables, Set to true in the task
#
# Summary of the edit (per the commit message / xml diff):
#  - FirewallPolicies_postrule (sheet1): remove the "PostRule" column, set the
#    Operation cell (A2) from "add" to "delete".
#  - NetworkAddresses / NetworkObjects / Services / ServiceGroups: move the
#    trailing "Operation" column to the front (column A) and change its value
#    from "add" to "delete".
#  - View-only housekeeping: selections / active tab / zoom tweaks to match
#    the saved workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) FirewallPolicies_postrule: drop the PostRule column (column S), and flip
#    the Operation value to "delete".
# ---------------------------------------------------------------------------
$wsPost = $wb.Worksheets.Item("FirewallPolicies_postrule")
$wsPost.Columns("S").Delete()
$wsPost.Range("A2").Value = "delete"

# ---------------------------------------------------------------------------
# 2) NetworkAddresses: move Operation (last column, E) to the front.
# ---------------------------------------------------------------------------
$wsAddr = $wb.Worksheets.Item("NetworkAddresses")
$wsAddr.Range("E1:E7").Cut()
$wsAddr.Columns("A").Insert()
$wsAddr.Range("A2:A7").Value = "delete"

# ---------------------------------------------------------------------------
# 3) NetworkObjects: move Operation (last column, E) to the front.
# ---------------------------------------------------------------------------
$wsObj = $wb.Worksheets.Item("NetworkObjects")
$wsObj.Range("E1:E5").Cut()
$wsObj.Columns("A").Insert()
$wsObj.Range("A2:A5").Value = "delete"

# ---------------------------------------------------------------------------
# 4) Services: move Operation (last column, F) to the front.
# ---------------------------------------------------------------------------
$wsSvc = $wb.Worksheets.Item("Services")
$wsSvc.Range("F1:F4").Cut()
$wsSvc.Columns("A").Insert()
$wsSvc.Range("A2:A4").Value = "delete"

# ---------------------------------------------------------------------------
# 5) ServiceGroups: move Operation (last column, E) to the front.
# ---------------------------------------------------------------------------
$wsGrp = $wb.Worksheets.Item("ServiceGroups")
$wsGrp.Range("E1:E3").Cut()
$wsGrp.Columns("A").Insert()
$wsGrp.Range("A2:A3").Value = "delete"

# ---------------------------------------------------------------------------
# 6) View/selection housekeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$wsPrerule = $wb.Worksheets.Item("FirewallPolicies_prerule")
$wsPrerule.Activate()
$wsPrerule.Range("E17").Select()

$wsAddr.Activate()
$wsAddr.Range("B12").Select()

$wsObj.Activate()
$wsObj.Range("C12").Select()

$wsSvc.Activate()
$excel.ActiveWindow.Zoom = 145
$wsSvc.Range("A8").Select()

$wsGrp.Activate()
$wsGrp.Range("A4").Select()

# FirewallPolicies_postrule becomes the active/selected tab last.
$wsPost.Activate()
$wsPost.Range("A3").Select()
